# Power comparisons for 3x2.
# Updates a set of numeric-looking text cells on the active sheet to reflect
# refreshed power-comparison figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = "0.2016 -0.0233"

    "D3"  = "< -0.0591"
    "E3"  = "< -0.0591"
    "F3"  = "< -0.0591"
    "G3"  = "0.6838 0.0152"
    "H3"  = "< -0.0591"
    "I3"  = "< -0.0591"
    "J3"  = "< -0.0591"
    "K3"  = "< -0.0591"
    "L3"  = "< -0.0591"
    "M3"  = "< -0.0591"
    "N3"  = "< -0.0591"
    "O3"  = "< -0.0591"

    "E4"  = "0.4463 0"
    "G4"  = "0.9802 0.0743"
    "J4"  = "0.6094 0"
    "K4"  = "0.4463 0"
    "L4"  = "0.6094 0"
    "M4"  = "0.4851 0"
    "N4"  = "0.4463 0"
    "O4"  = "0.4851 0"

    "F5"  = "0.5143 0"
    "G5"  = "0.9802 0.0743"
    "H5"  = "0.5143 0"
    "I5"  = "0.5143 0"
    "J5"  = "0.3206 0"
    "L5"  = "0.3206 0"

    "G6"  = "0.9802 0.0743"
    "J6"  = "0.6094 0"
    "K6"  = "0.4463 0"
    "L6"  = "0.6094 0"
    "M6"  = "0.4851 0"
    "N6"  = "0.4463 0"
    "O6"  = "0.4851 0"

    "H7"  = "< -0.0743"
    "I7"  = "< -0.0743"
    "J7"  = "< -0.0743"
    "K7"  = "< -0.0743"
    "L7"  = "< -0.0743"
    "M7"  = "< -0.0743"
    "N7"  = "< -0.0743"
    "O7"  = "< -0.0743"

    "J8"  = "0.6094 0"
    "K8"  = "0.4463 0"
    "L8"  = "0.6094 0"
    "M8"  = "0.4851 0"
    "N8"  = "0.4463 0"
    "O8"  = "0.4851 0"

    "J9"  = "0.6094 0"
    "K9"  = "0.4463 0"
    "L9"  = "0.6094 0"
    "M9"  = "0.4851 0"
    "N9"  = "0.4463 0"
    "O9"  = "0.4851 0"

    "K10" = "0.6592 0"
    "M10" = "0.681 0"
    "N10" = "0.6592 0"
    "O10" = "0.681 0"

    "L11" = "0.3206 0"

    "M12" = "0.681 0"
    "N12" = "0.6592 0"
    "O12" = "0.681 0"

    "N13" = "0.7089 0"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
